# Update countries & provincias Spain
#
# The underlying COVID country table (sheet "Pais") is kept sorted by
# "Casos totales" (column B, descending). This update refreshes several
# countries' figures, which also re-ranks (re-orders) a handful of rows
# that were tied/close in total cases, plus bumps the "updated at"
# timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a 19 de Mayo de 2020 a las 10:05" -> "...10:35"
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 19 de Mayo de 2020 a las 10:35'

# --- Row 36: Indonesia - refreshed totals
$ws.Cells.Item(36, 2).Value = 18496
$ws.Cells.Item(36, 3).Value = 486
$ws.Cells.Item(36, 4).Value = 4467
$ws.Cells.Item(36, 5).Value = 12808
$ws.Cells.Item(36, 7).Value = 30
$ws.Cells.Item(36, 8).Value = 1221

# --- Rows 44-46: Filipinas jumps ahead of Egipto / Republica Dominicana
# Row 44 becomes Filipinas with its newly refreshed totals
$ws.Cells.Item(44, 1).Value = 'Filipinas'
$ws.Cells.Item(44, 2).Value = 12942
$ws.Cells.Item(44, 3).Value = 224
$ws.Cells.Item(44, 4).Value = 2843
$ws.Cells.Item(44, 5).Value = 9262
$ws.Cells.Item(44, 7).Value = 6
$ws.Cells.Item(44, 8).Value = 837

# Row 45 becomes Egipto (its previous, row-44 figures)
$ws.Cells.Item(45, 1).Value = 'Egipto'
$ws.Cells.Item(45, 2).Value = 12764
$ws.Cells.Item(45, 4).Value = 3440
$ws.Cells.Item(45, 5).Value = 8679
$ws.Cells.Item(45, 8).Value = 645

# Row 46 becomes Republica Dominicana (its previous, row-45 figures)
$ws.Cells.Item(46, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(46, 2).Value = 12725
$ws.Cells.Item(46, 4).Value = 6613
$ws.Cells.Item(46, 5).Value = 5678
$ws.Cells.Item(46, 8).Value = 434

# Row 47 (Corea del Sur) is unchanged, row 48 Dinamarca refreshed on its own:
$ws.Cells.Item(48, 2).Value = 11044
$ws.Cells.Item(48, 3).Value = 76
$ws.Cells.Item(48, 5).Value = 1195

# --- Row 65: Oman - refreshed totals
$ws.Cells.Item(65, 2).Value = 5671
$ws.Cells.Item(65, 3).Value = 292
$ws.Cells.Item(65, 4).Value = 1574
$ws.Cells.Item(65, 5).Value = 4071
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 26

# --- Rows 196-197: Belice overtakes Nueva Caledonia
$ws.Cells.Item(196, 1).Value = 'Belice'
$ws.Cells.Item(196, 4).Value = 16
$ws.Cells.Item(196, 8).Value = 2

$ws.Cells.Item(197, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(197, 4).Value = 18
$ws.Cells.Item(197, 8).Value = 0

# --- Rows 209-211: Groenlandia & Montserrat overtake Seychelles
$ws.Cells.Item(209, 1).Value = 'Groenlandia'

$ws.Cells.Item(210, 1).Value = 'Montserrat'
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = 'Seychelles'
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# --- Rows 215-216: Bonaire, San Eustaquio y Saba overtakes San Bartolome
$ws.Cells.Item(215, 1).Value = 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(216, 1).Value = 'San Bartolome'
